$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "nvOfx812"
$ws.Range("B2").Value = 23080714
$ws.Range("C2").Value = "gkicyka73"
$ws.Range("D2").Value = "Z&#We6f3"
$ws.Range("F2").Value = "eToTjsaw"
$ws.Range("G2").Value = "wxeQ"
